$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write the new date label as a formula literal first (so Excel doesn't
# auto-convert the "dd-mm-yyyy" looking text into a real date/serial
# number), then flatten it down to a plain value via copy/paste-special.
# This keeps the resulting cell a plain shared-string cell with no
# extra/incidental number-format styling, matching how the rest of the
# column (A2:A14) is stored.
$ws.Range("A15").Formula = '="02-11-2021"'
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B15").Value = 810
$ws.Range("C15").Value = 825
